# Sync automático del tracker - 2025-08-31 20:21:00 UTC
# Updates rows 49-60 (matches played on 2025-08-30) from "Pending" to their
# final "Completed" results: Result, Resultado_Real, Profit, ROI, Enviado.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2025-08-31 19:17:12"

$updates = @(
    @{ Row = 49; Result = "Draw";     Real = "Fallo";   Profit = -2.3;  Roi = -100 },
    @{ Row = 50; Result = "Home Win"; Real = "Acierto"; Profit = 1.72;  Roi = 115 },
    @{ Row = 51; Result = "Away Win"; Real = "Acierto"; Profit = 3.02;  Roi = 57 },
    @{ Row = 52; Result = "Home Win"; Real = "Fallo";   Profit = -2.1;  Roi = -100 },
    @{ Row = 53; Result = "Home Win"; Real = "Acierto"; Profit = 2.92;  Roi = 55 },
    @{ Row = 54; Result = "Home Win"; Real = "Acierto"; Profit = 2.85;  Roi = 75 },
    @{ Row = 55; Result = "Home Win"; Real = "Acierto"; Profit = 2.12;  Roi = 40 },
    @{ Row = 56; Result = "Draw";     Real = "Fallo";   Profit = -2.7;  Roi = -100 },
    @{ Row = 57; Result = "Away Win"; Real = "Acierto"; Profit = 3.24;  Roi = 83 },
    @{ Row = 58; Result = "Draw";     Real = "Fallo";   Profit = -1;    Roi = -100 },
    @{ Row = 59; Result = "Home Win"; Real = "Acierto"; Profit = 2.7;   Roi = 55 },
    @{ Row = 60; Result = "Away Win"; Real = "Acierto"; Profit = 1.16;  Roi = 145 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 12).Value = "Completed"   # L: Status
    $ws.Cells.Item($r, 13).Value = $u.Result     # M: Result
    $ws.Cells.Item($r, 14).Value = $u.Real       # N: Resultado_Real
    $ws.Cells.Item($r, 15).Value = $u.Profit     # O: Profit
    $ws.Cells.Item($r, 16).Value = $u.Roi        # P: ROI
    $ws.Cells.Item($r, 17).Value = $timestamp    # Q: Enviado
}
